$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number need to be forced to
# text (Excel would otherwise silently coerce "504.86" -> numeric 504.86,
# with float rounding, and drop the original text formatting like "5.00").
$textCells = @(
    'D5',
    'D6',
    'D10',
    'D13',
    'D19',
    'D21',
    'D22',
    'D23',
    'D25',
    'D26',
    'D28',
    'D30',
    'D31',
    'D34',
    'D35',
    'D36',
    'D38',
    'D39',
    'D40',
    'D41',
    'D44',
    'D45',
    'D46',
    'D47',
    'D48'
)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '57.268.63'
$ws.Range('E2').Value = '  -0.92%  '
$ws.Range('D3').Value = '2.399.40'
$ws.Range('E3').Value = '  -1.91%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '504.86'
$ws.Range('E5').Value = '  -3.64%  '
$ws.Range('D6').Value = '132.41'
$ws.Range('E6').Value = '  +1.94%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  -0.56%  '
$ws.Range('D9').Value = '2.426.63'
$ws.Range('E9').Value = '  -0.87%  '
$ws.Range('D10').Value = '0.0973'
$ws.Range('E10').Value = '  -0.28%  '
$ws.Range('E11').Value = '  -0.90%  '
$ws.Range('E12').Value = '  +0.44%  '
$ws.Range('D13').Value = '4.66'
$ws.Range('E13').Value = '  -5.20%  '
$ws.Range('D14').Value = '2.833.13'
$ws.Range('E14').Value = '  -1.65%  '
$ws.Range('D15').Value = '57.087.32'
$ws.Range('E15').Value = '  -1.08%  '
$ws.Range('E16').Value = '  +0.75%  '
$ws.Range('E17').Value = '  +0.77%  '
$ws.Range('D18').Value = '2.421.74'
$ws.Range('E18').Value = '  -1.15%  '
$ws.Range('D19').Value = '10.26'
$ws.Range('E19').Value = '  -0.94%  '
$ws.Range('E20').Value = '  -0.27%  '
$ws.Range('D21').Value = '312.97'
$ws.Range('E21').Value = '  -0.52%  '
$ws.Range('D22').Value = '6.41'
$ws.Range('E22').Value = '  +5.17%  '
$ws.Range('D23').Value = '0.996'
$ws.Range('E23').Value = '  -0.29%  '
$ws.Range('E24').Value = '  -2.14%  '
$ws.Range('D25').Value = '65.28'
$ws.Range('E25').Value = '  +0.61%  '
$ws.Range('D26').Value = '0.993'
$ws.Range('E26').Value = '  -0.91%  '
$ws.Range('D27').Value = '2.508.95'
$ws.Range('E27').Value = '  -2.09%  '
$ws.Range('D28').Value = '0.384'
$ws.Range('E28').Value = '  -6.65%  '
$ws.Range('E29').Value = '  -2.80%  '
$ws.Range('D30').Value = '7.46'
$ws.Range('E30').Value = '  +3.47%  '
$ws.Range('D31').Value = '172.95'
$ws.Range('E31').Value = '  +0.61%  '
$ws.Range('D32').Value = '0.0₃0733'
$ws.Range('E32').Value = '  -0.01%  '
$ws.Range('E33').Value = '  -0.68%  '
$ws.Range('D34').Value = '6.20'
$ws.Range('E34').Value = '  +1.70%  '
$ws.Range('D35').Value = '1.15'
$ws.Range('E35').Value = '  +0.30%  '
$ws.Range('D36').Value = '0.998'
$ws.Range('E36').Value = '  -0.10%  '
$ws.Range('E37').Value = '  -0.69%  '
$ws.Range('D38').Value = '18.17'
$ws.Range('E38').Value = '  +2.02%  '
$ws.Range('D39').Value = '1.23'
$ws.Range('E39').Value = '  +4.17%  '
$ws.Range('D40').Value = '3.84'
$ws.Range('E40').Value = '  +0.89%  '
$ws.Range('D41').Value = '0.818'
$ws.Range('E41').Value = '  +3.66%  '
$ws.Range('E42').Value = '  -0.50%  '
$ws.Range('E43').Value = '  -0.34%  '
$ws.Range('D44').Value = '133.46'
$ws.Range('E44').Value = '  +7.80%  '
$ws.Range('D45').Value = '3.40'
$ws.Range('E45').Value = '  +0.08%  '
$ws.Range('D46').Value = '5.00'
$ws.Range('E46').Value = '  +4.11%  '
$ws.Range('D47').Value = '259.13'
$ws.Range('E47').Value = '  -2.01%  '
$ws.Range('D48').Value = '0.566'
$ws.Range('E48').Value = '  -2.90%  '
$ws.Range('E49').Value = '  -1.52%  '
$ws.Range('E50').Value = '  +0.49%  '
$ws.Range('E51').Value = '  +0.89%  '

# Restore the default "Normal" style on the forced-text cells so their
# style index matches the original (unstyled) cells.
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).Style = "Normal"
}
